# Sara-Alert-Format-Invalid-Monitorees.xlsx
# Adds three new "Race" related columns to the end of the header row on Sheet1:
#   CV1 = "Race Unknown"
#   CW1 = "Race Other"
#   CX1 = "Race Refused to Answer"
# The new header cells pick up the same formatting as the preceding header
# cell (CU1), and the active selection is moved to reflect the user having
# just finished typing the new headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Copy the formatting of the last existing header cell (CU1, style index 7)
# onto the three new header cells so they match the rest of the header row.
$lastHeader = $ws.Range("CU1")
$newHeaders = $ws.Range("CV1:CX1")
$lastHeader.Copy()
$newHeaders.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new header text.
$ws.Range("CV1").Value = "Race Unknown"
$ws.Range("CW1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"

# Scroll / select to roughly match where the editor ended up after adding
# the new columns.
$win = $excel.ActiveWindow
$win.ScrollColumn = 93
$ws.Range("CX6").Select()
